$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the data range A3:C102 by column A (ascending), matching the
# author's re-sort of the "Unmatched Publishers" table by publisher name
# instead of by the numeric match-count column B.
$sortRange = $ws.Range("A3:C102")
$key1 = $ws.Range("A3:A102")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($key1)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4142   # xlNo
$ws.Sort.Apply()

# The resort leaves the previously-blank row 72 with no data, so Excel's
# used range shrinks by one row once that gap is gone.
$ws.Range("A49").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1

$wb.Windows.Item(1).WindowState = -4143 # xlNormal, no-op safeguard

# Restore the original selection and move the window's top-left cell back
# to A1 (so the saved view no longer pins row 49 at the top).
$ws.Range("A3:C101").Select()
$excel.ActiveWindow.ScrollRow = 1

$wb.Windows.Item(1).Left = 9320
